# Updated symbol list with GitHub Actions - refresh cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> { column letter -> new literal text value }
# Values must stay plain text (matching the source data feed format), so each
# cell is forced to Text format before the write, then the formatting footprint
# is cleared again so no stray number-format style sticks to the cell.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextCell 2 4 "303.27"
Set-TextCell 2 5 "-2.01%"

# Row 3
Set-TextCell 3 4 "35.87"
Set-TextCell 3 5 "0.87%"

# Row 4
Set-TextCell 4 4 "5.065"
Set-TextCell 4 5 "-1.15%"

# Row 5
Set-TextCell 5 4 "0.08060"
Set-TextCell 5 5 "-1.56%"

# Row 6
Set-TextCell 6 4 "1.934"
Set-TextCell 6 5 "-5.55%"

# Row 7
Set-TextCell 7 4 "7.788"
Set-TextCell 7 5 "-2.12%"

# Row 8
Set-TextCell 8 4 "0.9299"
Set-TextCell 8 5 "0.20%"

# Row 9
Set-TextCell 9 4 "0.1515"
Set-TextCell 9 5 "40.76%"

# Row 10
Set-TextCell 10 4 "0.1907"
Set-TextCell 10 5 "-0.56%"

# Row 11
Set-TextCell 11 4 "0.08987"
Set-TextCell 11 5 "-5.96%"

# Row 12
Set-TextCell 12 4 "0.03459"
Set-TextCell 12 5 "-4.47%"

# Row 13
Set-TextCell 13 4 "0.09842"
Set-TextCell 13 5 "-0.65%"

# Row 14
Set-TextCell 14 4 "0.001397"
Set-TextCell 14 5 "-2.23%"

# Row 15
Set-TextCell 15 4 "0.005745"
Set-TextCell 15 5 "0.88%"

# Row 16
Set-TextCell 16 4 "3.539"
Set-TextCell 16 5 "1.97%"

# Row 17
Set-TextCell 17 4 "4.054"
Set-TextCell 17 5 "-1.92%"

# Row 18
Set-TextCell 18 5 "2.99%"

# Row 19
Set-TextCell 19 4 "0.3444"
Set-TextCell 19 5 "0.85%"

# Row 20
Set-TextCell 20 5 "-0.97%"

# Row 21
Set-TextCell 21 4 "5.037"
Set-TextCell 21 5 "-1.34%"

# Row 22
Set-TextCell 22 4 "0.2394"
Set-TextCell 22 5 "9.04%"

# Row 23
Set-TextCell 23 4 "0.04499"
Set-TextCell 23 5 "-1.13%"

# Row 24
Set-TextCell 24 5 "-1.51%"

# Row 25
Set-TextCell 25 4 "0.004808"
Set-TextCell 25 5 "0.63%"

# Row 26
Set-TextCell 26 5 "-2.11%"

# Row 27
Set-TextCell 27 5 "-32.39%"

# Row 39
Set-TextCell 39 4 "0.01876"
Set-TextCell 39 5 "-4.40%"

# Row 40
Set-TextCell 40 4 "0.04795"
Set-TextCell 40 5 "-1.69%"

# Row 41
Set-TextCell 41 5 "7.40%"

# Row 42
Set-TextCell 42 4 "0.007284"
Set-TextCell 42 5 "-5.25%"

# Row 43
Set-TextCell 43 4 "0.1344"
Set-TextCell 43 5 "-2.56%"

# Row 44
Set-TextCell 44 4 "0.002104"
Set-TextCell 44 5 "-0.74%"

# Row 45
Set-TextCell 45 4 "0.009721"
Set-TextCell 45 5 "-15.81%"

# Row 46
Set-TextCell 46 4 "0.00006221"
Set-TextCell 46 5 "-4.53%"

# Row 47
Set-TextCell 47 5 "-0.52%"

# Row 48
Set-TextCell 48 5 "0.48%"

# Row 50
Set-TextCell 50 5 "-0.52%"

# Row 51
Set-TextCell 51 5 "-0.52%"
